# For PO Form Modify
#
# The "Terms and Conditions" block on the PO form used to contain plain
# static labels. This swaps them for merge-field-ready labels (matching the
# supplier/date/etc. fields used elsewhere on the form) and leaves the sheet
# scrolled/selected on the block that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "General Condition: #GeneralCondition"
$ws.Range("A22").Value = "Quality Management System Requirements: #POQMSRequirement`n"
$ws.Range("A23").Value = "Quality: #POQuality"
$ws.Range("A24").Value = "Packing & Forwarding: #POPackForward"
$ws.Range("A25").Value = "Mode of Payment: #ModeOfPayment"
$ws.Range("A26").Value = "Payment Terms: #PaymentTerms"
$ws.Range("A27").Value = "Mode of Transport: #ModeOfTransport"
$ws.Range("A28").Value = "Any Other Requirements: #AnyOtherRequirements"
$ws.Range("A29").Value = "PO Validity  : #POValidity"

# Scroll the sheet so row 18 is at the top and select A29, matching the
# saved view state of the edited workbook.
$excel.Goto($ws.Range("A18"), $true)
[void]$ws.Range("A29").Select()
